# Edit slide 8 ("Team World Wide" deck) - replace the parking-meter analysis
# paragraphs in the content placeholder with the updated commentary about
# pricing / EV charging spots.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)          # "Content Placeholder 2"
$tr = $sh.TextFrame.TextRange

# Start from a clean slate so the new text doesn't inherit stray run
# boundaries from whatever used to occupy a given paragraph index.
$tr.Text = ""

$para1 = "Using this, we may be able to exponentialize pricing or reduce max parking times to create a better flow."
$para2 = "We also wanted to analyze which part of town needs more parking spaces based on the availability of the current parking meters."
$para3Prefix = "Using this, we may be able to find which part of town requires electric "
$para3Suffix = "charging spots."

# `\r` is how PowerPoint's TextRange represents paragraph breaks. Add one
# extra trailing break so we can cleanly trim the body back down to
# exactly three paragraphs afterwards.
$tr.Text = $para1 + "`r" + $para2 + "`r" + $para3Prefix + $para3Suffix + "`r"

# Remove the trailing empty paragraph introduced by the final `\r`.
$tr.Paragraphs(4, 1).Delete()

# Split the third paragraph into two runs so the "charging spots." tail
# carries its own run (matching the target formatting split).
$para3 = $tr.Paragraphs(3, 1)
$tail = $para3.Characters($para3Prefix.Length + 1, $para3.Length - $para3Prefix.Length)
$tail.Text = $tail.Text
